$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-4 down to 3-5)
$ws.Rows.Item(2).Insert()

# Remove any formatting copied onto the inserted row from the row above
$ws.Range("A2").ClearFormats()

# Set the value of the newly inserted cell
$ws.Range("A2").Value = "uniquely_henrietta"

# Update the selection to C3 as per the diff
$ws.Range("C3").Select()
